$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.876.04"
$ws.Range("E2").Value = "  +6.96%  "
$ws.Range("D3").Value = "3.019.55"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.75"
$ws.Range("E5").Value = "  +3.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.43"
$ws.Range("E6").Value = "  +9.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "3.014.69"
$ws.Range("E8").Value = "  +4.28%  "
$ws.Range("E9").Value = "  +2.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.99"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.156"
$ws.Range("E11").Value = "  +6.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("E12").Value = "  +5.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000249"
$ws.Range("E13").Value = "  +8.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.50"
$ws.Range("E14").Value = "  +8.34%  "
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "65.953.99"
$ws.Range("E16").Value = "  +7.09%  "
$ws.Range("D17").Value = "3.521.12"
$ws.Range("E17").Value = "  +4.42%  "
$ws.Range("E18").Value = "  +6.87%  "
$ws.Range("D19").Value = "3.020.03"
$ws.Range("E19").Value = "  +4.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "465.48"
$ws.Range("E20").Value = "  +7.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.85"
$ws.Range("E21").Value = "  +6.21%  "
$ws.Range("E22").Value = "  +4.49%  "
$ws.Range("E23").Value = "  +8.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.20"
$ws.Range("E24").Value = "  +3.60%  "
$ws.Range("E25").Value = "  +5.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("E26").Value = "  +11.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.81"
$ws.Range("E27").Value = "  +9.33%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.96"
$ws.Range("E29").Value = "  +13.90%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.42"
$ws.Range("E30").Value = "  +18.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0000106"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.62"
$ws.Range("E32").Value = "  +5.01%  "
$ws.Range("E33").Value = "  +5.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.10"
$ws.Range("E34").Value = "  +6.13%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +4.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.82"
$ws.Range("E37").Value = "  +8.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.19"
$ws.Range("E38").Value = "  +13.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.07"
$ws.Range("E39").Value = "  +8.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.31"
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "45.00"
$ws.Range("E41").Value = "  +13.45%  "
$ws.Range("E42").Value = "  +7.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.300"
$ws.Range("E43").Value = "  +12.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.52"
$ws.Range("E44").Value = "  +3.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "391.91"
$ws.Range("E45").Value = "  +13.38%  "
$ws.Range("D46").Value = "2.801.46"
$ws.Range("E46").Value = "  +4.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0356"
$ws.Range("E47").Value = "  +6.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.38"
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.78"
$ws.Range("E50").Value = "  +10.24%  "
$ws.Range("E51").Value = "  +4.44%  "
